$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Continue the hours log on Sheet1 with three more entries (rows 34-36).
# Copy the date-formatted cell above so the new date cells inherit the
# same number format (m/d/yyyy) as the rest of column A.
$ws.Range("A33").Copy() | Out-Null
$ws.Range("A34:A36").PasteSpecial(-4122) | Out-Null

# Row 34: 2010-03-23, 2.5 hours, Group Meeting
$ws.Range("A34").Value = "3/23/2010"
$ws.Range("B34").Value = 2.5
$ws.Range("C34").Value = "Group Meeting"

# Row 35: 2010-03-23, 1 hour, Weekly Meeting
$ws.Range("A35").Value = "3/23/2010"
$ws.Range("B35").Value = 1
$ws.Range("C35").Value = "Weekly Meeting"

# Row 36: 2010-03-26, 0.5 hours, Skype Meeting
$ws.Range("A36").Value = "3/26/2010"
$ws.Range("B36").Value = 0.5
$ws.Range("C36").Value = "Skype Meeting"

# Keep the selection/active cell consistent with the appended data
$ws.Range("A37").Select()
